$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assign RA names for each schedule date per the updated roster order
# (sorted by date added: Rachel, unassigned, Sadaf, Jamie, Zoya)
$ws.Cells.Item(4, 2).Value = "Sadaf"
$ws.Cells.Item(5, 2).Value = "Jamie"
$ws.Cells.Item(6, 2).Value = "Zoya"
$ws.Cells.Item(7, 2).Value = "Rachel"
$ws.Cells.Item(8, 2).Value = "Jamie"
$ws.Cells.Item(9, 2).Value = "Rachel"
$ws.Cells.Item(11, 2).Value = "Sadaf"
$ws.Cells.Item(12, 2).Value = "Rachel"
$ws.Cells.Item(13, 2).Value = "Zoya"
$ws.Cells.Item(14, 2).Value = "Jamie"
$ws.Cells.Item(15, 2).Value = "Rachel"
$ws.Cells.Item(18, 2).Value = "Rachel"
$ws.Cells.Item(19, 2).Value = "Jamie"
$ws.Cells.Item(20, 2).Value = "Zoya"
$ws.Cells.Item(21, 2).Value = "Zoya"
$ws.Cells.Item(28, 2).Value = "Sadaf"
$ws.Cells.Item(30, 2).Value = "Rachel"
$ws.Cells.Item(32, 2).Value = "Rachel"
$ws.Cells.Item(39, 2).Value = "Zoya"
$ws.Cells.Item(41, 2).Value = "Jamie"
$ws.Cells.Item(42, 2).Value = "Rachel"
$ws.Cells.Item(43, 2).Value = "unassigned"
$ws.Cells.Item(45, 2).Value = "Jamie"
$ws.Cells.Item(46, 2).Value = "Zoya"
$ws.Cells.Item(47, 2).Value = "Zoya"
$ws.Cells.Item(48, 2).Value = "Jamie"
$ws.Cells.Item(49, 2).Value = "Sadaf"
$ws.Cells.Item(50, 2).Value = "Zoya"
$ws.Cells.Item(51, 2).Value = "Rachel"
$ws.Cells.Item(53, 2).Value = "Zoya"
$ws.Cells.Item(55, 2).Value = "Rachel"
$ws.Cells.Item(56, 2).Value = "Jamie"
$ws.Cells.Item(57, 2).Value = "Rachel"
$ws.Cells.Item(58, 2).Value = "Jamie"
$ws.Cells.Item(59, 2).Value = "Rachel"
$ws.Cells.Item(60, 2).Value = "Zoya"
$ws.Cells.Item(63, 2).Value = "Jamie"
$ws.Cells.Item(64, 2).Value = "Sadaf"
$ws.Cells.Item(65, 2).Value = "Jamie"
$ws.Cells.Item(66, 2).Value = "Zoya"
$ws.Cells.Item(67, 2).Value = "Zoya"
$ws.Cells.Item(70, 2).Value = "Sadaf"
$ws.Cells.Item(71, 2).Value = "unassigned"
$ws.Cells.Item(72, 2).Value = "Zoya"
$ws.Cells.Item(73, 2).Value = "Jamie"
$ws.Cells.Item(74, 2).Value = "unassigned"
$ws.Cells.Item(75, 2).Value = "Zoya"
$ws.Cells.Item(76, 2).Value = "Sadaf"
$ws.Cells.Item(77, 2).Value = "unassigned"
$ws.Cells.Item(78, 2).Value = "unassigned"

# Widen column A to fit the (now longer-looking) schedule column
$ws.Range("A1").EntireColumn.ColumnWidth = 33.3
